# Add two new header/data columns (I, F-suffixed "I0"/"IF") to the sheet,
# matching the existing header styling and plain numeric data cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1) — same formatted look as the existing headers.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, border, center/top alignment) from the
# existing "IP" header cell (H1) onto the two new header cells so they
# share the same cell style as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data cells (row 2) — plain numeric values, unstyled like the other
# data cells in that row.
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
